# Apply updates to Jogos_do_Dia_Betfair_Back_Lay_2025-10-06.xlsx
# - 180 odds cells updated on existing rows 2-14
# - 2 new match rows appended (15, 16) with full data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update existing odds values (rows 2-14) ---
$ws.Range("F2").Value = 2.64
$ws.Range("I2").Value = 2.7
$ws.Range("N2").Value = 5.4
$ws.Range("O2").Value = 1.19
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 1.63
$ws.Range("S2").Value = 2.42
$ws.Range("T2").Value = 1.53
$ws.Range("V2").Value = 1.59
$ws.Range("X2").Value = 28
$ws.Range("Y2").Value = 980
$ws.Range("AL2").Value = 32
$ws.Range("I3").Value = 8.8
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.55
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 1.88
$ws.Range("Q3").Value = 1.94
$ws.Range("T3").Value = 2.06
$ws.Range("U3").Value = 1.79
$ws.Range("V3").Value = 1.13
$ws.Range("W3").Value = 2.72
$ws.Range("X3").Value = 18
$ws.Range("Y3").Value = 28
$ws.Range("Z3").Value = 85
$ws.Range("AB3").Value = 8.8
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 38
$ws.Range("AF3").Value = 10.5
$ws.Range("AG3").Value = 12.5
$ws.Range("AH3").Value = 32
$ws.Range("AJ3").Value = 17
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 55
$ws.Range("S4").Value = 3.25
$ws.Range("Y4").Value = 13
$ws.Range("AB4").Value = 16
$ws.Range("AC4").Value = 9.8
$ws.Range("AD4").Value = 14
$ws.Range("AG4").Value = 17
$ws.Range("AJ4").Value = 70
$ws.Range("G5").Value = 1.42
$ws.Range("H5").Value = 5.9
$ws.Range("J5").Value = 4
$ws.Range("L5").Value = 1.27
$ws.Range("N5").Value = 4.1
$ws.Range("R5").Value = 1.52
$ws.Range("S5").Value = 2.3
$ws.Range("T5").Value = 1.94
$ws.Range("U5").Value = 1.83
$ws.Range("X5").Value = 28
$ws.Range("Y5").Value = 42
$ws.Range("AB5").Value = 11.5
$ws.Range("AC5").Value = 16
$ws.Range("AD5").Value = 46
$ws.Range("AF5").Value = 10.5
$ws.Range("AG5").Value = 13
$ws.Range("AH5").Value = 36
$ws.Range("AJ5").Value = 14
$ws.Range("AK5").Value = 17.5
$ws.Range("AL5").Value = 44
$ws.Range("AN5").Value = 6.4
$ws.Range("F6").Value = 3.9
$ws.Range("G6").Value = 4.8
$ws.Range("I6").Value = 1.94
$ws.Range("J6").Value = 4.2
$ws.Range("P6").Value = 2.62
$ws.Range("Q6").Value = 1.45
$ws.Range("R6").Value = 1.66
$ws.Range("S6").Value = 2.3
$ws.Range("V6").Value = 2.06
$ws.Range("W6").Value = 1.26
$ws.Range("X6").Value = 34
$ws.Range("Y6").Value = 17
$ws.Range("Z6").Value = 18
$ws.Range("AA6").Value = 26
$ws.Range("AB6").Value = 28
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 13.5
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 44
$ws.Range("AG6").Value = 22
$ws.Range("AH6").Value = 19.5
$ws.Range("AI6").Value = 32
$ws.Range("AJ6").Value = 95
$ws.Range("AK6").Value = 50
$ws.Range("AL6").Value = 50
$ws.Range("AM6").Value = 75
$ws.Range("AN6").Value = 36
$ws.Range("AO6").Value = 9.4
$ws.Range("F7").Value = 2.28
$ws.Range("I7").Value = 3.9
$ws.Range("L7").Value = 1.35
$ws.Range("P7").Value = 1.82
$ws.Range("Q7").Value = 1.97
$ws.Range("V7").Value = 1.34
$ws.Range("X7").Value = 980
$ws.Range("Y7").Value = 980
$ws.Range("Z7").Value = 980
$ws.Range("AB7").Value = 12
$ws.Range("AC7").Value = 9.4
$ws.Range("AD7").Value = 980
$ws.Range("AE7").Value = 980
$ws.Range("AF7").Value = 980
$ws.Range("AG7").Value = 980
$ws.Range("AH7").Value = 980
$ws.Range("AI7").Value = 65
$ws.Range("AJ7").Value = 980
$ws.Range("AK7").Value = 980
$ws.Range("AL7").Value = 980
$ws.Range("AN7").Value = 980
$ws.Range("AO7").Value = 980
$ws.Range("P9").Value = 1.4
$ws.Range("I10").Value = 2.88
$ws.Range("L10").Value = 1.44
$ws.Range("X10").Value = 16
$ws.Range("Y10").Value = 14
$ws.Range("Z10").Value = 22
$ws.Range("AA10").Value = 50
$ws.Range("AB10").Value = 14
$ws.Range("AC10").Value = 9.6
$ws.Range("AD10").Value = 15.5
$ws.Range("AE10").Value = 38
$ws.Range("AF10").Value = 23
$ws.Range("AG10").Value = 15.5
$ws.Range("AH10").Value = 22
$ws.Range("AI10").Value = 55
$ws.Range("AJ10").Value = 55
$ws.Range("AK10").Value = 38
$ws.Range("AL10").Value = 55
$ws.Range("AN10").Value = 34
$ws.Range("AO10").Value = 32
$ws.Range("H11").Value = 1.82
$ws.Range("Q11").Value = 2.34
$ws.Range("T11").Value = 2.22
$ws.Range("U11").Value = 1.64
$ws.Range("I12").Value = 2.58
$ws.Range("L12").Value = 1.73
$ws.Range("M12").Value = 1.13
$ws.Range("S12").Value = 6.8
$ws.Range("T12").Value = 2.5
$ws.Range("U12").Value = 1.57
$ws.Range("V12").Value = 1.63
$ws.Range("W12").Value = 1.35
$ws.Range("X12").Value = 8.4
$ws.Range("AD12").Value = 16.5
$ws.Range("AK12").Value = 100
$ws.Range("AO12").Value = 70
$ws.Range("G13").Value = 2.26
$ws.Range("I13").Value = 5.8
$ws.Range("J13").Value = 3.05
$ws.Range("K13").Value = 3.4
$ws.Range("L13").Value = 1.48
$ws.Range("M13").Value = 1.11
$ws.Range("N13").Value = 2.5
$ws.Range("O13").Value = 1.54
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.38
$ws.Range("T13").Value = 2.16
$ws.Range("U13").Value = 1.7
$ws.Range("V13").Value = 1.21
$ws.Range("W13").Value = 1.79
$ws.Range("X13").Value = 10.5
$ws.Range("F14").Value = 1.48
$ws.Range("G14").Value = 1.59
$ws.Range("H14").Value = 6.8
$ws.Range("I14").Value = 11.5
$ws.Range("J14").Value = 3.5
$ws.Range("K14").Value = 5.4
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 2.84
$ws.Range("O14").Value = 1.37
$ws.Range("P14").Value = 1.74
$ws.Range("Q14").Value = 2.08
$ws.Range("R14").Value = 1.27
$ws.Range("S14").Value = 3.45
$ws.Range("T14").Value = 2.02
$ws.Range("U14").Value = 1.68
$ws.Range("V14").Value = 1.09
$ws.Range("W14").Value = 2.64

# --- 2) Append two new match rows (15 and 16) ---
# Date column (B) needs a quote-prefix so Excel keeps it as literal text
# instead of auto-converting "2025-10-06" into a date serial number.

# Row 15
$ws.Range("A15").Value = 'Argentinian Primera Division'
$ws.Range("B15").Value = '''2025-10-06'
$ws.Range("C15").Value = '21:00:00'
$ws.Range("D15").Value = 'Racing Club'
$ws.Range("E15").Value = 'Independiente Rivadavia'
$ws.Range("F15").Value = 1.75
$ws.Range("G15").Value = 1.91
$ws.Range("H15").Value = 5
$ws.Range("I15").Value = 6.8
$ws.Range("J15").Value = 3.55
$ws.Range("K15").Value = 4
$ws.Range("L15").Value = 1.48
$ws.Range("M15").Value = 1.08
$ws.Range("N15").Value = 3.1
$ws.Range("O15").Value = 1.41
$ws.Range("P15").Value = 1.69
$ws.Range("Q15").Value = 2.2
$ws.Range("R15").Value = 1.26
$ws.Range("S15").Value = 4
$ws.Range("T15").Value = 2.02
$ws.Range("U15").Value = 1.78
$ws.Range("V15").Value = 1.2
$ws.Range("W15").Value = 2.08
$ws.Range("X15").Value = 13.5
$ws.Range("Y15").Value = 19
$ws.Range("Z15").Value = 50
$ws.Range("AA15").Value = 180
$ws.Range("AB15").Value = 9
$ws.Range("AC15").Value = 10
$ws.Range("AD15").Value = 27
$ws.Range("AE15").Value = 110
$ws.Range("AF15").Value = 12
$ws.Range("AG15").Value = 13
$ws.Range("AH15").Value = 29
$ws.Range("AI15").Value = 120
$ws.Range("AJ15").Value = 21
$ws.Range("AK15").Value = 27
$ws.Range("AL15").Value = 60
$ws.Range("AM15").Value = 210
$ws.Range("AN15").Value = 18.5
$ws.Range("AO15").Value = 150

# Row 16
$ws.Range("A16").Value = 'Colombian Primera A'
$ws.Range("B16").Value = '''2025-10-06'
$ws.Range("C16").Value = '21:30:00'
$ws.Range("D16").Value = 'La Equidad'
$ws.Range("E16").Value = 'Once Caldas'
$ws.Range("F16").Value = 3.25
$ws.Range("G16").Value = 3.85
$ws.Range("H16").Value = 2.34
$ws.Range("I16").Value = 2.6
$ws.Range("J16").Value = 3.05
$ws.Range("K16").Value = 3.5
$ws.Range("L16").Value = 1.42
$ws.Range("M16").Value = 1.08
$ws.Range("N16").Value = 2.86
$ws.Range("O16").Value = 1.44
$ws.Range("P16").Value = 1.63
$ws.Range("Q16").Value = 2.26
$ws.Range("R16").Value = 1.23
$ws.Range("S16").Value = 3.9
$ws.Range("T16").Value = 1.9
$ws.Range("U16").Value = 1.89
$ws.Range("V16").Value = 1.62
$ws.Range("W16").Value = 1.35
$ws.Range("X16").Value = 980
$ws.Range("Y16").Value = 10.5
$ws.Range("Z16").Value = 980
$ws.Range("AA16").Value = 980
$ws.Range("AB16").Value = 980
$ws.Range("AC16").Value = 8.8
$ws.Range("AD16").Value = 980
$ws.Range("AE16").Value = 980
$ws.Range("AF16").Value = 980
$ws.Range("AG16").Value = 980
$ws.Range("AH16").Value = 980
$ws.Range("AI16").Value = 65
$ws.Range("AJ16").Value = 85
$ws.Range("AK16").Value = 60
$ws.Range("AL16").Value = 80
$ws.Range("AM16").Value = 1000
$ws.Range("AN16").Value = 75
$ws.Range("AO16").Value = 980
